$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.194.46'
$ws.Range("E2").Value = '  +1.14%  '

# Row 3
$ws.Range("D3").Value = '1.834.81'
$ws.Range("E3").Value = '  +1.05%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.012'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +1.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.25%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.010'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.97%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4709'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.54%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3685'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.26%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07417'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.55%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8820'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.31%  '

# Row 12
$ws.Range("D12").Value = '1.830.23'
$ws.Range("E12").Value = '  +3.05%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07333'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.63%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.473'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.98%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.86'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.67%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.567'
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.014'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.18%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008791'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.79%  '

# Row 19
$ws.Range("E19").Value = '  +0.96%  '

# Row 20
$ws.Range("E20").Value = '  +0.42%  '

# Row 21
$ws.Range("D21").Value = '27.221.97'
$ws.Range("E21").Value = '  +1.17%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.309'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.47%  '

# Row 23
$ws.Range("E23").Value = '  +1.45%  '

# Row 24
$ws.Range("D24").Value = '2.055.74'
$ws.Range("E24").Value = '  +1.64%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.905'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.65%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.53'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.43%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.99%  '

# Row 28
$ws.Range("E28").Value = '  -1.25%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.279'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.45%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.77'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.99%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08926'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.05%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7595'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.94%  '

# Row 33
$ws.Range("E33").Value = '  +0.92%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.545'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.41%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.945'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.84%  '

# Row 36
$ws.Range("E36").Value = '  +1.02%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.105'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.71%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05340'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.55%  '

# Row 39
$ws.Range("E39").Value = '  +0.12%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.007'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.27%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.327'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.13%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.411'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.80%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5351'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.21%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1663'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.06%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.540'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.10%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4954'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.26%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.54'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.03%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.011'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.07%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.673'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.13%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '104.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.19%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06325'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.67%  '
